# Updates the yearly income-statement figures on the "Overview" sheet
# (rows 11-27, columns D:H) with the newly reported financial data,
# and zeroes out the two rows that previously held a placeholder "-".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$data = @{
    11 = @(6206795, 10130258, 18044120, 23112181, 37692415)
    12 = @(-3948341, -6511883, -12697090, -17756829, -26317310)
    13 = @(2258454, 3618375, 5347030, 5355352, 11375105)
    14 = @(-603615, -721899, -1348877, -2260746, -3335380)
    15 = @(0, 0, 0, 0, 0)
    16 = @(11635, 2315, 190125, 62065, -42234)
    17 = @(1666474, 2898791, 4188278, 3156671, 7997491)
    18 = @(-347155, -636364, -954588, -1490615, -2491810)
    19 = @(199820, 411980, 1090767, 1182099, 1152674)
    20 = @(1519139, 2674407, 4324457, 2848155, 6658355)
    21 = @(-136992, -144540, -155529, -128161, -52040)
    22 = @(1382147, 2529867, 4168928, 2719994, 6606315)
    23 = @(0, 0, 0, 0, 0)
    24 = @(1382147, 2529867, 4168928, 2719994, 6606315)
    25 = @(2304, 2530, 4169, 1360, 1652)
    26 = @(600000, 1000000, 1000000, 2000000, 4000000)
    27 = @(346, 632, 1042, 680, 1652)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("D$row").Value = $values[0]
    $ws.Range("E$row").Value = $values[1]
    $ws.Range("F$row").Value = $values[2]
    $ws.Range("G$row").Value = $values[3]
    $ws.Range("H$row").Value = $values[4]
}
